$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 166, pushing existing rows 166-204 down to 167-205.
$ws.Rows.Item(166).Insert()

# Populate the newly inserted row 166 with the new data record.
$ws.Cells.Item(166, 1).Value = 8
$ws.Cells.Item(166, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(166, 3).Value = "Coquimbo"
$ws.Cells.Item(166, 4).Value = 44508
$ws.Cells.Item(166, 5).Value = 4
$ws.Cells.Item(166, 6).Value = 100112032
$ws.Cells.Item(166, 7).Value = "Zapallo italiano"
$ws.Cells.Item(166, 8).Value = "Sin especificar"
$ws.Cells.Item(166, 9).Value = "Primera"
$ws.Cells.Item(166, 10).Value = 560
$ws.Cells.Item(166, 11).Value = 10000
$ws.Cells.Item(166, 12).Value = 11000
$ws.Cells.Item(166, 13).Value = 10500
$ws.Cells.Item(166, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(166, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(166, 16).Value = 150
$ws.Cells.Item(166, 17).Value = 70
$ws.Cells.Item(166, 18).Value = "Hortaliza"
